# Update the "Förändrad" (changed) date column (C) for rows 2-23
# from serial date 45183 (2023-09-14) to 45184 (2023-09-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
